$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1:AF1 - copy the existing header formatting (bold,
# bordered, centered) from AC1 so the new headers match the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2..45: team record (Wins/Losses/Ties), same for every player row
for ($r = 2; $r -le 45; $r++) {
    $ws.Cells.Item($r, 30).Value = 89   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 73   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
